# Fix the typo "informative" -> "informatif" in the "Thématique du site"
# paragraph, then drop Word's usual "_GoBack" bookmark at the spot of the
# last edit (immediately after the corrected word), matching what Word
# itself records when you make and save an edit like this.

$d = $word.ActiveDocument

$r = $d.Content
$found = $r.Find.Execute("informative", $true, $false, $false, $false, $false, `
                          $true, 1, $false, "", 0)

if ($found -and $r.Find.Found) {
    $r.Text = "informatif"

    $editEnd = $r.End
    $goBack = $d.Range($editEnd, $editEnd)
    $d.Bookmarks.Add("_GoBack", $goBack)
}
